# Applies numeric-value corrections to the Chocobo_Profits workbook sheets
# (currentAveragePrice / LevePrice / LeveProfit columns H-N) per the commit diff.
$wb = $excel.ActiveWorkbook

# ALC row 19: Unbreak My Heart
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19,8).Value = 833546.5
$ws.Cells.Item(19,9).Value = 1904941.9
$ws.Cells.Item(19,11).Value = 1904941.9
$ws.Cells.Item(19,13).Value = -1904766.9

# ALC row 116: Growing Up
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116,8).Value = 565978.75
$ws.Cells.Item(116,9).Value = 1002780.5
$ws.Cells.Item(116,10).Value = 19976.5
$ws.Cells.Item(116,11).Value = 1002780.5
$ws.Cells.Item(116,12).Value = 19976.5
$ws.Cells.Item(116,13).Value = -999338.5
$ws.Cells.Item(116,14).Value = -26860.5

# ALC row 123: Nearly Bare
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(123,8).Value = 42980
$ws.Cells.Item(123,10).Value = 42980
$ws.Cells.Item(123,12).Value = 42980
$ws.Cells.Item(123,14).Value = -52780

# ALC row 132: Fast-forwarding Flora
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132,8).Value = 85410.78999999999
$ws.Cells.Item(132,9).Value = 101625.95
$ws.Cells.Item(132,11).Value = 304877.85
$ws.Cells.Item(132,13).Value = -302347.85

# ALC row 137: Cutting Edge of Culinary Quality
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137,8).Value = 2870.7273
$ws.Cells.Item(137,9).Value = 1993.36
$ws.Cells.Item(137,10).Value = 5612.5
$ws.Cells.Item(137,11).Value = 5980.08
$ws.Cells.Item(137,12).Value = 16837.5
$ws.Cells.Item(137,13).Value = -3430.08
$ws.Cells.Item(137,14).Value = -21937.5

# ALC row 138: All-night Crafting
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138,8).Value = 2766.0476
$ws.Cells.Item(138,9).Value = 2048.9285
$ws.Cells.Item(138,10).Value = 2909.4714
$ws.Cells.Item(138,11).Value = 6146.7855
$ws.Cells.Item(138,12).Value = 8728.414199999999
$ws.Cells.Item(138,13).Value = -1006.7855
$ws.Cells.Item(138,14).Value = -19008.4142

# ARM row 2: Ain't Got No Ingots
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2,8).Value = 882.5
$ws.Cells.Item(2,9).Value = 965
$ws.Cells.Item(2,10).Value = 800
$ws.Cells.Item(2,11).Value = 965
$ws.Cells.Item(2,12).Value = 800
$ws.Cells.Item(2,13).Value = -852
$ws.Cells.Item(2,14).Value = -1026

# ARM row 32: Ingot We Trust
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32,8).Value = 12013.412
$ws.Cells.Item(32,9).Value = 8111.59
$ws.Cells.Item(32,10).Value = 17260.69
$ws.Cells.Item(32,11).Value = 8111.59
$ws.Cells.Item(32,12).Value = 17260.69
$ws.Cells.Item(32,13).Value = -7824.59
$ws.Cells.Item(32,14).Value = -17834.69

# ARM row 74: As the Bolt Flies
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74,8).Value = 1640.881
$ws.Cells.Item(74,9).Value = 1164.9445
$ws.Cells.Item(74,11).Value = 1164.9445
$ws.Cells.Item(74,13).Value = -290.9445000000001

# ARM row 77: Heavy Metal Banned (L)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77,8).Value = 1640.881
$ws.Cells.Item(77,9).Value = 1164.9445
$ws.Cells.Item(77,11).Value = 5824.7225
$ws.Cells.Item(77,13).Value = -1456.7225

# ARM row 88: The Mast Chance
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(88,8).Value = 5558097.5
$ws.Cells.Item(88,9).Value = 16668141
$ws.Cells.Item(88,10).Value = 3075.75
$ws.Cells.Item(88,11).Value = 16668141
$ws.Cells.Item(88,12).Value = 3075.75
$ws.Cells.Item(88,13).Value = -16667735
$ws.Cells.Item(88,14).Value = -3887.75

# ARM row 91: The Rose and the Riveter (L)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(91,8).Value = 5558097.5
$ws.Cells.Item(91,9).Value = 16668141
$ws.Cells.Item(91,10).Value = 3075.75
$ws.Cells.Item(91,11).Value = 16668141
$ws.Cells.Item(91,12).Value = 3075.75
$ws.Cells.Item(91,13).Value = -16666737
$ws.Cells.Item(91,14).Value = -5883.75

# ARM row 116: No Scope
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116,8).Value = 882.5
$ws.Cells.Item(116,9).Value = 965
$ws.Cells.Item(116,10).Value = 800
$ws.Cells.Item(116,11).Value = 965
$ws.Cells.Item(116,12).Value = 800
$ws.Cells.Item(116,13).Value = 1329
$ws.Cells.Item(116,14).Value = -5388

# ARM row 124: Ace of Gloves
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(124,8).Value = 50000
$ws.Cells.Item(124,10).Value = 50000
$ws.Cells.Item(124,12).Value = 50000
$ws.Cells.Item(124,14).Value = -59820

# ARM row 132: Don't Bore Me, Ore Me
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132,8).Value = 3114.1035
$ws.Cells.Item(132,9).Value = 1441.2307
$ws.Cells.Item(132,10).Value = 4473.3125
$ws.Cells.Item(132,11).Value = 4323.6921
$ws.Cells.Item(132,12).Value = 13419.9375
$ws.Cells.Item(132,13).Value = -1793.6921
$ws.Cells.Item(132,14).Value = -18479.9375

# ARM row 138: Don't Ask about the Rivets
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(138,8).Value = 70000
$ws.Cells.Item(138,10).Value = 70000
$ws.Cells.Item(138,12).Value = 70000
$ws.Cells.Item(138,14).Value = -80280

# BSM row 3: Hells Bells
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3,8).Value = 882.5
$ws.Cells.Item(3,9).Value = 965
$ws.Cells.Item(3,10).Value = 800
$ws.Cells.Item(3,11).Value = 965
$ws.Cells.Item(3,12).Value = 800
$ws.Cells.Item(3,13).Value = -851
$ws.Cells.Item(3,14).Value = -1028

# BSM row 22: Riveting Run
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22,8).Value = 527.5714
$ws.Cells.Item(22,9).Value = 198.4
$ws.Cells.Item(22,10).Value = 1350.5
$ws.Cells.Item(22,11).Value = 198.4
$ws.Cells.Item(22,12).Value = 1350.5
$ws.Cells.Item(22,13).Value = -25.40000000000001
$ws.Cells.Item(22,14).Value = -1696.5

# BSM row 26: Unseamly Conditions
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(26,8).Value = 0
$ws.Cells.Item(26,9).Value = 0
$ws.Cells.Item(26,11).Value = 0
$ws.Cells.Item(26,13).ClearContents()

# BSM row 86: Through Thick and Thin
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86,8).Value = 2041
$ws.Cells.Item(86,10).Value = 2990
$ws.Cells.Item(86,12).Value = 2990
$ws.Cells.Item(86,14).Value = -5236

# BSM row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89,8).Value = 2041
$ws.Cells.Item(89,10).Value = 2990
$ws.Cells.Item(89,12).Value = 14950
$ws.Cells.Item(89,14).Value = -26182

# BSM row 94: High Steal
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94,8).Value = 2768.2856
$ws.Cells.Item(94,10).Value = 2000
$ws.Cells.Item(94,12).Value = 2000
$ws.Cells.Item(94,14).Value = -2902

# BSM row 103: The Bigger the Blade
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(103,8).Value = 0
$ws.Cells.Item(103,10).Value = 0
$ws.Cells.Item(103,12).Value = 0
$ws.Cells.Item(103,14).ClearContents()

# BSM row 105: Ingot to Wing It
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105,8).Value = 2654.4546
$ws.Cells.Item(105,9).Value = 2569.9
$ws.Cells.Item(105,10).Value = 3500
$ws.Cells.Item(105,11).Value = 2569.9
$ws.Cells.Item(105,12).Value = 3500
$ws.Cells.Item(105,13).Value = -822.9000000000001
$ws.Cells.Item(105,14).Value = -6994

# BSM row 137: Dagger Swagger
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(137,8).Value = 33086.6
$ws.Cells.Item(137,10).Value = 33086.6
$ws.Cells.Item(137,12).Value = 33086.6
$ws.Cells.Item(137,14).Value = -43286.6

# CRP row 132: Hull Lotta Damage
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132,8).Value = 2964.4517
$ws.Cells.Item(132,9).Value = 1536
$ws.Cells.Item(132,11).Value = 4608
$ws.Cells.Item(132,13).Value = -2078

# GSM row 15: The Tusk at Hand
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(15,8).Value = 42556.855
$ws.Cells.Item(15,10).Value = 42556.855
$ws.Cells.Item(15,12).Value = 42556.855
$ws.Cells.Item(15,14).Value = -43132.855

# GSM row 62: The Goggles, They Do Naught
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(62,8).Value = 0
$ws.Cells.Item(62,10).Value = 0
$ws.Cells.Item(62,12).Value = 0
$ws.Cells.Item(62,14).ClearContents()

# GSM row 63: Not on My Table
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(63,8).Value = 19900
$ws.Cells.Item(63,10).Value = 19900
$ws.Cells.Item(63,12).Value = 19900
$ws.Cells.Item(63,14).Value = -21272

# GSM row 64: Halonic Hermeneutics
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(64,8).Value = 50000
$ws.Cells.Item(64,10).Value = 50000
$ws.Cells.Item(64,12).Value = 50000
$ws.Cells.Item(64,14).Value = -50496

# GSM row 65: Peril Never Wore Safety Goggles (L)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(65,8).Value = 0
$ws.Cells.Item(65,10).Value = 0
$ws.Cells.Item(65,12).Value = 0
$ws.Cells.Item(65,14).ClearContents()

# GSM row 66: Heinz's Dilemma (L)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(66,8).Value = 19900
$ws.Cells.Item(66,10).Value = 19900
$ws.Cells.Item(66,12).Value = 59700
$ws.Cells.Item(66,14).Value = -66564

# GSM row 67: Transposing Theology (L)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(67,8).Value = 50000
$ws.Cells.Item(67,10).Value = 50000
$ws.Cells.Item(67,12).Value = 50000
$ws.Cells.Item(67,14).Value = -51716

# GSM row 74: The Unfortunate Retirony
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(74,8).Value = 35582.75
$ws.Cells.Item(74,10).Value = 35582.75
$ws.Cells.Item(74,12).Value = 35582.75
$ws.Cells.Item(74,14).Value = -37454.75

# GSM row 77: Life Ends at Retirement (L)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(77,8).Value = 35582.75
$ws.Cells.Item(77,10).Value = 35582.75
$ws.Cells.Item(77,12).Value = 106748.25
$ws.Cells.Item(77,14).Value = -116108.25

# GSM row 80: Needs More Prayerbell
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80,8).Value = 125001500

# GSM row 81: The Grander Temple
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(81,8).Value = 42556.855
$ws.Cells.Item(81,10).Value = 42556.855
$ws.Cells.Item(81,12).Value = 42556.855
$ws.Cells.Item(81,14).Value = -44552.855

# GSM row 83: With a Noise That Reaches Heaven (L)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83,8).Value = 125001500

# GSM row 84: Man with a Dragon Earring (L)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(84,8).Value = 42556.855
$ws.Cells.Item(84,10).Value = 42556.855
$ws.Cells.Item(84,12).Value = 127670.565
$ws.Cells.Item(84,14).Value = -137654.565

# GSM row 123: Workplace Workout
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(123,8).Value = 15261.883
$ws.Cells.Item(123,10).Value = 15261.883
$ws.Cells.Item(123,12).Value = 15261.883
$ws.Cells.Item(123,14).Value = -20161.883

# LTW row 22: Skin off Their Backs
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22,8).Value = 2450.2273
$ws.Cells.Item(22,9).Value = 2244.6667
$ws.Cells.Item(22,11).Value = 2244.6667
$ws.Cells.Item(22,13).Value = -1949.6667

# LTW row 27: Fire and Hide
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27,8).Value = 2450.2273
$ws.Cells.Item(27,9).Value = 2244.6667
$ws.Cells.Item(27,11).Value = 2244.6667
$ws.Cells.Item(27,13).Value = -2137.6667

# LTW row 46: Supply Side Logic
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46,8).Value = 2892.8572
$ws.Cells.Item(46,9).Value = 5000
$ws.Cells.Item(46,10).Value = 2541.6667
$ws.Cells.Item(46,11).Value = 5000
$ws.Cells.Item(46,12).Value = 2541.6667
$ws.Cells.Item(46,13).Value = -4812
$ws.Cells.Item(46,14).Value = -2917.6667

# LTW row 93: Hide to Go Seek
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93,8).Value = 10102672
$ws.Cells.Item(93,9).Value = 13889874
$ws.Cells.Item(93,10).Value = 3466.6667
$ws.Cells.Item(93,11).Value = 13889874
$ws.Cells.Item(93,12).Value = 3466.6667
$ws.Cells.Item(93,13).Value = -13888626
$ws.Cells.Item(93,14).Value = -5962.6667

# LTW row 134: Freezing Fingers
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(134,8).Value = 48719.168
$ws.Cells.Item(134,10).Value = 48719.168
$ws.Cells.Item(134,12).Value = 48719.168
$ws.Cells.Item(134,14).Value = -58859.168

# LTW row 136: Respect for Br'aax
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136,8).Value = 3176.7273
$ws.Cells.Item(136,9).Value = 1669.6111
$ws.Cells.Item(136,10).Value = 4985.2666
$ws.Cells.Item(136,11).Value = 5008.8333
$ws.Cells.Item(136,12).Value = 14955.7998
$ws.Cells.Item(136,13).Value = -2458.8333
$ws.Cells.Item(136,14).Value = -20055.7998

# LTW row 138: Freezing Toes
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(138,8).Value = 118899
$ws.Cells.Item(138,10).Value = 118899
$ws.Cells.Item(138,12).Value = 118899
$ws.Cells.Item(138,14).Value = -129179

# WVR row 54: No Country for Cold Men
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(54,8).Value = 17315.375
$ws.Cells.Item(54,10).Value = 17315.375
$ws.Cells.Item(54,12).Value = 17315.375
$ws.Cells.Item(54,14).Value = -18355.375

# WVR row 113: A Tender Table
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113,8).Value = 8098.846
$ws.Cells.Item(113,9).Value = 11496.223
$ws.Cells.Item(113,10).Value = 454.75
$ws.Cells.Item(113,11).Value = 34488.669
$ws.Cells.Item(113,12).Value = 1364.25
$ws.Cells.Item(113,13).Value = -32318.669
$ws.Cells.Item(113,14).Value = -5704.25

# WVR row 130: Skill Cap
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(130,8).Value = 40055.91
$ws.Cells.Item(130,10).Value = 40055.91
$ws.Cells.Item(130,12).Value = 40055.91
$ws.Cells.Item(130,14).Value = -50095.91

# WVR row 132: Comfy Cabins
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132,8).Value = 10106222
$ws.Cells.Item(132,9).Value = 7985.2856
$ws.Cells.Item(132,11).Value = 23955.8568
$ws.Cells.Item(132,13).Value = -21425.8568
